$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A currently only holds the header "Nombre" in A1 (rest of column A is
# empty); column B holds the real header "Correo" in B1 and the email
# addresses in B2:B132. Deleting the entire column A removes the "Nombre"
# header and shifts column B (Correo + all emails) one column to the left,
# landing on column A.
$ws.Range("A1").EntireColumn.Delete()

# Update the selection to match the post-edit state.
$ws.Range("D9").Select()
